# Added Week 15 simulations
$wb = $excel.ActiveWorkbook

# OFF sheet - row 3 updates
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 171
$wsOff.Range("C3").Value = 118
$wsOff.Range("D3").Value = 45
$wsOff.Range("F3").Value = 2

# DEF sheet - row 3 updates
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 158
$wsDef.Range("C3").Value = 116
$wsDef.Range("D3").Value = 30
$wsDef.Range("E3").Value = 14
